$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 0.3430696666666667
$ws.Range("H2").Value = 1.029209
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 153.5290173333333
$ws.Range("N2").Value = 460.587052
$ws.Range("O2").Value = 0.3172206968818489
$ws.Range("P2").Value = 0.317220696881849
$ws.Range("Q2").Value = 52.67114880020755
$ws.Range("R2").Value = 474.040339201868
$ws.Range("S2").Value = 0.3172206968818489
$ws.Range("T2").Value = 0.317220696881849

# Row 3
$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 0.3430696666666667
$ws.Range("H3").Value = 1.029209
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 168.7997026666667
$ws.Range("N3").Value = 506.3991080000001
$ws.Range("O3").Value = 0.3487728915577651
$ws.Range("P3").Value = 0.3487728915577651
$ws.Range("Q3").Value = 57.91005772728578
$ws.Range("R3").Value = 521.1905195455721
$ws.Range("S3").Value = 0.3487728915577651
$ws.Range("T3").Value = 0.3487728915577651

# Row 4
$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 0.3430696666666667
$ws.Range("H4").Value = 1.029209
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 68.09032333333333
$ws.Range("N4").Value = 204.27097
$ws.Range("O4").Value = 0.1406878008722904
$ws.Range("P4").Value = 0.1406878008722904
$ws.Range("Q4").Value = 23.35972452919222
$ws.Range("R4").Value = 210.23752076273
$ws.Range("S4").Value = 0.1406878008722904
$ws.Range("T4").Value = 0.1406878008722904

# Row 5
$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 0.3430696666666667
$ws.Range("H5").Value = 1.029209
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 93.562673
$ws.Range("N5").Value = 280.688019
$ws.Range("O5").Value = 0.1933186106880956
$ws.Range("P5").Value = 0.1933186106880956
$ws.Range("Q5").Value = 32.09851503855234
$ws.Range("R5").Value = 288.886635346971
$ws.Range("S5").Value = 0.1933186106880956
$ws.Range("T5").Value = 0.1933186106880956
